# VG-FS-ADO-Sync.xlsx edit:
#  - Add a new "FS-Field-Type" column (E) to the SingleField sheet, with
#    text/date type markers for each FS field row.
#  - Remove the devops_status / System.State mapping row (old row 10).
#  - Update sheet views so SingleField becomes the active tab (instead of
#    ProductsData), matching the single-cell selections left behind.

$wb = $excel.ActiveWorkbook

$singleField = $wb.Worksheets.Item("SingleField")
$productsData = $wb.Worksheets.Item("ProductsData")

# --- SingleField (sheet1) ---------------------------------------------

# Remove the last row (devops_status -> System.State, ADO_TO_FS) entirely.
$singleField.Rows.Item(10).Delete()

# New header cell for column E, copying the bold/bordered header style
# used by the other header cells in row 1.
$singleField.Cells.Item(1, 4).Copy($singleField.Cells.Item(1, 5))
$singleField.Cells.Item(1, 5).Value = "FS-Field-Type"

# Per-row FS field type values (blank for id/department_id rows).
$singleField.Cells.Item(3, 5).Value = "text"
$singleField.Cells.Item(4, 5).Value = "date"
$singleField.Cells.Item(5, 5).Value = "date"
$singleField.Cells.Item(6, 5).Value = "text"
$singleField.Cells.Item(7, 5).Value = "text"
$singleField.Cells.Item(9, 5).Value = "date"

# Widen the new column to fit its header/content.
$singleField.Columns.Item(5).ColumnWidth = 18.14

# --- ProductsData (sheet6) ---------------------------------------------

# Collapse its multi-cell selection down to a single cell (no longer the
# active tab) before switching focus away from it.
$productsData.Range("H2").Select()

# --- Switch the active tab back to SingleField --------------------------
$singleField.Range("E2").Select()
